$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column Q: year 2020, mirroring the format of column P (2019) ---
$ws.Range("P4:P5").Copy()
$ws.Range("Q4:Q5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(4, 17).Value = 2020
$ws.Cells.Item(5, 17).Value = 16.7

# --- Row 5 data cells switch number format from "0.00" to "0.0" ---
$ws.Range("D5:Q5").NumberFormat = "0.0"

# --- View: scroll so column C is left-most visible, select N12:N13 ---
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("N12:N13").Select()
